$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.738.57'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '2.245.24'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.83'
$ws.Range('E5').Value = '  +2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.53'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  -1.61%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.37'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0828'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.69'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').Value = '2.588.24'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.858'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('D17').Value = '2.249.24'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('D18').Value = '43.692.96'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.80'
$ws.Range('E19').Value = '  -4.34%  '
$ws.Range('D20').Value = '0.0₃0984'
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.48'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.20'
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.17'
$ws.Range('E23').Value = '  +0.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '235.92'
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.16'
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.09'
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').Value = '  -1.78%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.08'
$ws.Range('E29').Value = '  +7.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.29'
$ws.Range('E30').Value = '  -1.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '160.54'
$ws.Range('E31').Value = '  +3.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.16'
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0852'
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.70'
$ws.Range('E34').Value = '  -2.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.19'
$ws.Range('E35').Value = '  +0.63%  '
$ws.Range('E36').Value = '  +9.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.94'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.119'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.80'
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.30'
$ws.Range('E40').Value = '  -3.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.53'
$ws.Range('E41').Value = '  +20.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0318'
$ws.Range('E42').Value = '  -1.87%  '
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').Value = '1.818.61'
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.73'
$ws.Range('E46').Value = '  +6.13%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '82.82'
$ws.Range('E47').Value = '  -4.88%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.23'
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '74.58'
$ws.Range('E49').Value = '  -4.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '58.78'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '103.70'
$ws.Range('E51').Value = '  +0.02%  '
